# Add the site URL to the workbook: a new "Sheet2" (placed after "Sheet1")
# holding a "url" label and the automation-practice site link, plus leave
# the selection on Sheet1 parked at E25 (matches the author's last-saved
# cursor position).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 so tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Write the URL value first, then the "url" label, so the shared-string
# table picks up the two new strings in that order (index 4 = link text,
# index 5 = "url").
$ws2.Range("A2").Value = "http://automationpractice.com/index.php"
$ws2.Range("A1").Value = "url"

# Restore the active sheet/selection to Sheet1 at E25.
$ws1.Range("E25").Select()
